# Trade #8 closed at 2026-02-16 22:58:03 - base_strategy UP +0.000%
# Append the new trade row (row 9) to both the "All Trades" and
# "base_strategy" worksheets.

$wb = $excel.ActiveWorkbook

$sheetNames = @("All Trades", "base_strategy")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    $row = 9

    # --- Numeric / plain-text columns (no auto-conversion risk) ---
    $ws.Cells.Item($row, 1).Value = 8                      # A: Trade #
    $ws.Cells.Item($row, 4).Value = "base_strategy"        # D: Strategy
    $ws.Cells.Item($row, 5).Value = "UP"                   # E: Side
    $ws.Cells.Item($row, 6).Value = 0.5                    # F: Entry Price
    $ws.Cells.Item($row, 8).Value = "OPEN"                 # H: Status
    $ws.Cells.Item($row, 9).Value = 0                      # I: P&L %
    $ws.Cells.Item($row, 10).Value = 0                     # J: P&L $
    $ws.Cells.Item($row, 11).Value = 100                   # K: Capital After
    $ws.Cells.Item($row, 12).Value = 0                     # L: Entry Slippage (bps)
    $ws.Cells.Item($row, 13).Value = 0                     # M: Exit Slippage (bps)
    $ws.Cells.Item($row, 14).Value = 0.6                   # N: Confidence
    $ws.Cells.Item($row, 15).Value = "Normal spread capture: 19600 bps" # O: Entry Reason
    $ws.Cells.Item($row, 17).Value = 0                     # Q: Duration (min)

    # --- Column B holds a date-shaped string ("2026-02-16") that Excel's
    # auto-detection would otherwise coerce into a real date serial value.
    # Force literal text via a leading quote-prefix, then strip the
    # resulting quotePrefix style back to Normal so the cell keeps plain
    # default formatting. ---
    $ws.Cells.Item($row, 2).Value = "'2026-02-16"
    $ws.Cells.Item($row, 2).Style = "Normal"

    # Column C ("22:58:03") is not auto-converted, but set the same way
    # for consistency/safety.
    $ws.Cells.Item($row, 3).Value = "22:58:03"

    # --- Columns G and P are empty-string text cells (Exit Price / Exit
    # Reason are blank while the trade is still OPEN). A plain "" assigns
    # a truly blank cell, but the source data stores an empty *text*
    # cell, so use the quote-prefix trick to get an empty string value,
    # then reset the style back to Normal. ---
    $ws.Cells.Item($row, 7).Value = "'"
    $ws.Cells.Item($row, 7).Style = "Normal"

    $ws.Cells.Item($row, 16).Value = "'"
    $ws.Cells.Item($row, 16).Style = "Normal"
}

Write-Output "Trade #8 appended to All Trades and base_strategy sheets"
